$d = $word.ActiveDocument

# The paragraph that holds the _GoBack bookmark is the last paragraph in the
# document body (Paragraphs.Count gives its 1-based position reliably; the
# Paragraph.Index property does not match this collection's numbering, so it
# is avoided here).
$bookmarkPos = $d.Paragraphs.Count

# Insert a new paragraph right before the bookmark paragraph and fill it with
# the "4.Tạo controller" heading line.
$d.Paragraphs.Item($bookmarkPos).Range.InsertParagraphBefore()
$d.Paragraphs.Item($bookmarkPos).Range.Text = "4.Tạo controller"

# The bookmark paragraph shifted down by one; insert the second new
# paragraph (the artisan command line) right before it as well.
$bookmarkPos = $d.Paragraphs.Count
$d.Paragraphs.Item($bookmarkPos).Range.InsertParagraphBefore()
$d.Paragraphs.Item($bookmarkPos).Range.Text = "php artisan make:controller Name_Table"

# Recompute the bookmark paragraph again and append a single space run after
# the existing bookmarkStart/bookmarkEnd markers (InsertAfter on the
# collapsed end-of-range keeps it after those markers).
$bookmarkPos = $d.Paragraphs.Count
$bmRange = $d.Paragraphs.Item($bookmarkPos).Range
$endRange = $d.Range($bmRange.End, $bmRange.End)
$endRange.InsertAfter(" ")

# Finally, add a trailing empty paragraph after the bookmark paragraph (at the
# very end of the document body, before the sectPr).
$bookmarkPos = $d.Paragraphs.Count
$d.Paragraphs.Item($bookmarkPos).Range.InsertParagraphAfter()
